{"js": "// The \"Function Summary\" table lists C functions alongside their Fortran\n// equivalents. The row for H5Ocopy had \"(none)\" listed as its Fortran\n// wrapper; the commit adds the newly available Fortran wrapper name\n// \"H5ocopy_f\" in its place.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Locate the \"Function Summary\" table: the one whose second row reads\n// \"C Function\" / \"Fortran\" in its first cell.\nlet targetTable = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  table.load(\"rowCount\");\n}\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  if (table.rowCount < 3) continue;\n  const headerCell = table.getCell(1, 0);\n  headerCell.load(\"value\");\n  await context.sync();\n  if (headerCell.value.indexOf(\"C Function\") !== -1 && headerCell.value.indexOf(\"Fortran\") !== -1) {\n    targetTable = table;\n    break;\n  }\n}\n\nif (!targetTable) {\n  throw new Error(\"Could not find the Function Summary table\");\n}\n\n// Within that table, find the row whose first paragraph is \"H5Ocopy\" and\n// replace the second paragraph (\"(none)\") with \"H5ocopy_f\".\ntargetTable.load(\"rowCount\");\nawait context.sync();\n\nlet updated = false;\nfor (let r = 0; r < targetTable.rowCount && !updated; r++) {\n  const cell = targetTable.getCell(r, 0);\n  const cellBody = cell.body;\n  cellBody.load(\"paragraphs/items/text\");\n  await context.sync();\n\n  const paras = cellBody.paragraphs.items;\n  if (paras.length >= 2 && paras[0].text.trim() === \"H5Ocopy\" && paras[1].text.trim() === \"(none)\") {\n    paras[1].insertText(\"H5ocopy_f\", Word.InsertLocation.replace);\n    await context.sync();\n    updated = true;\n  }\n}\n\nif (!updated) {\n  throw new Error(\"Could not find the H5Ocopy / (none) row to update\");\n}\n", "ps1": "# The \"Function Summary\" table lists C functions alongside their Fortran\n# equivalents. The row for H5Ocopy had \"(none)\" listed as its Fortran\n# wrapper; this commit adds the newly available Fortran wrapper name\n# \"H5ocopy_f\" in its place.\n\n$d = $word.ActiveDocument\n\n$targetTable = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $candidate = $d.Tables.Item($i)\n    if ($candidate.Rows.Count -ge 3) {\n        $headerText = $candidate.Cell(2, 1).Range.Text\n        if ($headerText -like \"*C Function*\" -and $headerText -like \"*Fortran*\") {\n            $targetTable = $candidate\n            break\n        }\n    }\n}\n\nif ($targetTable -eq $null) {\n    throw \"Could not find the Function Summary table\"\n}\n\n$updated = $false\nfor ($r = 1; $r -le $targetTable.Rows.Count; $r++) {\n    $cell = $targetTable.Cell($r, 1)\n    $cellText = $cell.Range.Text -replace \"`r`a\", \"|\" -replace \"`r\", \"|\"\n    $parts = $cellText.Split(\"|\")\n    if ($parts.Length -ge 2 -and $parts[0].Trim() -eq \"H5Ocopy\" -and $parts[1].Trim() -eq \"(none)\") {\n        $found = $cell.Range.Find.Execute(\"(none)\", $false, $false, $false, $false, $false, $true, 0, $false, \"H5ocopy_f\", 1)\n        if ($found) {\n            $updated = $true\n        }\n        break\n    }\n}\n\nif (-not $updated) {\n    throw \"Could not find the H5Ocopy / (none) row to update\"\n}\n"}
